$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 100004040
$ws.Range("I111").Value = 125004930
$ws.Range("J111").Value = 495
$ws.Range("K111").Value = 375014790
$ws.Range("L111").Value = 1485
$ws.Range("M111").Value = -375011723
$ws.Range("N111").Value = -7619
$ws.Range("H112").Value = 1770.3334
$ws.Range("I112").Value = 199.5
$ws.Range("J112").Value = 2084.5
$ws.Range("K112").Value = 598.5
$ws.Range("L112").Value = 6253.5
$ws.Range("M112").Value = 509.5
$ws.Range("N112").Value = -8469.5
$ws.Range("H132").Value = 2466.1143
$ws.Range("I132").Value = 2614.3704
$ws.Range("J132").Value = 1965.75
$ws.Range("K132").Value = 7843.111199999999
$ws.Range("L132").Value = 5897.25
$ws.Range("M132").Value = -5313.111199999999
$ws.Range("N132").Value = -10957.25
$ws.Range("H135").Value = 409.07144
$ws.Range("I135").Value = 363.76923
$ws.Range("K135").Value = 3273.92307
$ws.Range("M135").Value = -738.9230699999998
$ws.Range("H137").Value = 20919476
$ws.Range("I137").Value = 31252212
$ws.Range("J137").Value = 254000.75
$ws.Range("K137").Value = 93756636
$ws.Range("L137").Value = 762002.25
$ws.Range("M137").Value = -93754086
$ws.Range("N137").Value = -767102.25
$ws.Range("H138").Value = 3664.724
$ws.Range("I138").Value = 844.2174
$ws.Range("J138").Value = 5518.2
$ws.Range("K138").Value = 2532.6522
$ws.Range("L138").Value = 16554.6
$ws.Range("M138").Value = 2607.3478
$ws.Range("N138").Value = -26834.6

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1580.6538
$ws.Range("I2").Value = 1636.8422
$ws.Range("J2").Value = 1428.1428
$ws.Range("K2").Value = 1636.8422
$ws.Range("L2").Value = 1428.1428
$ws.Range("M2").Value = -1523.8422
$ws.Range("N2").Value = -1654.1428
$ws.Range("H61").Value = 2771.25
$ws.Range("I61").Value = 1025.5
$ws.Range("K61").Value = 1025.5
$ws.Range("M61").Value = -813.5
$ws.Range("H74").Value = 9686.154
$ws.Range("I74").Value = 12657.777
$ws.Range("K74").Value = 12657.777
$ws.Range("M74").Value = -11783.777
$ws.Range("H77").Value = 9686.154
$ws.Range("I77").Value = 12657.777
$ws.Range("K77").Value = 63288.885
$ws.Range("M77").Value = -58920.885
$ws.Range("H116").Value = 1580.6538
$ws.Range("I116").Value = 1636.8422
$ws.Range("J116").Value = 1428.1428
$ws.Range("K116").Value = 1636.8422
$ws.Range("L116").Value = 1428.1428
$ws.Range("M116").Value = 657.1578
$ws.Range("N116").Value = -6016.1428
$ws.Range("H123").Value = 41119.332
$ws.Range("J123").Value = 41119.332
$ws.Range("L123").Value = 41119.332
$ws.Range("N123").Value = -50919.332
$ws.Range("H132").Value = 2817.8057
$ws.Range("I132").Value = 2671.32
$ws.Range("J132").Value = 3150.7273
$ws.Range("K132").Value = 8013.960000000001
$ws.Range("L132").Value = 9452.1819
$ws.Range("M132").Value = -5483.960000000001
$ws.Range("N132").Value = -14512.1819
$ws.Range("H136").Value = 2771.25
$ws.Range("I136").Value = 1025.5
$ws.Range("K136").Value = 3076.5
$ws.Range("M136").Value = -526.5

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1580.6538
$ws.Range("I3").Value = 1636.8422
$ws.Range("J3").Value = 1428.1428
$ws.Range("K3").Value = 1636.8422
$ws.Range("L3").Value = 1428.1428
$ws.Range("M3").Value = -1522.8422
$ws.Range("N3").Value = -1656.1428
$ws.Range("H107").Value = 2749.2104
$ws.Range("I107").Value = 2866.7646
$ws.Range("J107").Value = 1750
$ws.Range("K107").Value = 2866.7646
$ws.Range("L107").Value = 1750
$ws.Range("M107").Value = -946.7646
$ws.Range("N107").Value = -5590

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2887.7222
$ws.Range("I16").Value = 1983
$ws.Range("J16").Value = 5240
$ws.Range("K16").Value = 1983
$ws.Range("L16").Value = 5240
$ws.Range("M16").Value = -1696
$ws.Range("N16").Value = -5814
$ws.Range("H74").Value = 25176.75
$ws.Range("I74").Value = 26400
$ws.Range("K74").Value = 26400
$ws.Range("M74").Value = -25526
$ws.Range("H77").Value = 25176.75
$ws.Range("I77").Value = 26400
$ws.Range("K77").Value = 79200
$ws.Range("M77").Value = -74832
$ws.Range("H108").Value = 29933.334
$ws.Range("J108").Value = 29933.334
$ws.Range("L108").Value = 29933.334
$ws.Range("N108").Value = -37613.334
$ws.Range("H113").Value = 2887.7222
$ws.Range("I113").Value = 1983
$ws.Range("J113").Value = 5240
$ws.Range("K113").Value = 1983
$ws.Range("L113").Value = 5240
$ws.Range("M113").Value = 187
$ws.Range("N113").Value = -9580
$ws.Range("H132").Value = 2115.2
$ws.Range("I132").Value = 529.125
$ws.Range("J132").Value = 3927.8572
$ws.Range("K132").Value = 1587.375
$ws.Range("L132").Value = 11783.5716
$ws.Range("M132").Value = 942.625
$ws.Range("N132").Value = -16843.5716
$ws.Range("H134").Value = 4275.1665
$ws.Range("I134").Value = 821.3333
$ws.Range("K134").Value = 2463.9999
$ws.Range("M134").Value = 71.0001000000002

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 732
$ws.Range("I131").Value = 382.17392
$ws.Range("J131").Value = 949.4595
$ws.Range("K131").Value = 1146.52176
$ws.Range("L131").Value = 2848.3785
$ws.Range("M131").Value = 3893.47824
$ws.Range("N131").Value = -12928.3785

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 6938.1904
$ws.Range("I113").Value = 1860.8889
$ws.Range("J113").Value = 10746.167
$ws.Range("K113").Value = 1860.8889
$ws.Range("L113").Value = 10746.167
$ws.Range("M113").Value = 309.1111000000001
$ws.Range("N113").Value = -15086.167
$ws.Range("H116").Value = 46979.25
$ws.Range("J116").Value = 46979.25
$ws.Range("L116").Value = 46979.25
$ws.Range("N116").Value = -56157.25

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5558.8887
$ws.Range("I61").Value = 6759.5
$ws.Range("J61").Value = 2128.5715
$ws.Range("K61").Value = 6759.5
$ws.Range("L61").Value = 2128.5715
$ws.Range("M61").Value = -6557.5
$ws.Range("N61").Value = -2532.5715
$ws.Range("H112").Value = 50000
$ws.Range("J112").Value = 50000
$ws.Range("L112").Value = 50000
$ws.Range("N112").Value = -52954
$ws.Range("H113").Value = 5558.8887
$ws.Range("I113").Value = 6759.5
$ws.Range("J113").Value = 2128.5715
$ws.Range("K113").Value = 6759.5
$ws.Range("L113").Value = 2128.5715
$ws.Range("M113").Value = -4589.5
$ws.Range("N113").Value = -6468.5715
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H128").Value = 40000
$ws.Range("J128").Value = 40000
$ws.Range("L128").Value = 40000
$ws.Range("N128").Value = -49960

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 809.3077
$ws.Range("I107").Value = 548.8570999999999
$ws.Range("J107").Value = 1113.1666
$ws.Range("K107").Value = 1646.5713
$ws.Range("L107").Value = 3339.4998
$ws.Range("M107").Value = 273.4287000000002
$ws.Range("N107").Value = -7179.4998
$ws.Range("H126").Value = 635.6957
$ws.Range("I126").Value = 654.5263
$ws.Range("K126").Value = 1963.5789
$ws.Range("M126").Value = 506.4211
$ws.Range("H132").Value = 1754.5186
$ws.Range("I132").Value = 1292.75
$ws.Range("J132").Value = 2426.182
$ws.Range("K132").Value = 3878.25
$ws.Range("L132").Value = 7278.545999999999
$ws.Range("M132").Value = -1348.25
$ws.Range("N132").Value = -12338.546
